# Generate Report for handback
#
# For each localized-language sheet ("zh-cn", "de-de") this:
#   - flips the Status column (B) from "Ready for handoff" to
#     "Handed back: in sync with en-US" for the two tracked source files
#   - fills in the (previously empty) "Latest Target File" (E) and
#     "Latest Handback File" (F) columns with the same md / xlf file names
#     already used for the handoff columns (A / C), including hyperlinks
#   - stamps "Latest Handback DateTime" (G) with the handback timestamp

$wb = $excel.ActiveWorkbook

function Get-LinkAddress($sheet, $cellRef) {
    foreach ($h in $sheet.Hyperlinks) {
        if ($h.Range.Address() -eq $cellRef) {
            return $h.Address()
        }
    }
    return $null
}

function Set-HandbackRow($sheet, $row, $handbackDatetime) {
    $aRef = "`$A`$" + $row
    $cRef = "`$C`$" + $row

    $aAddr = Get-LinkAddress $sheet $aRef
    $cAddr = Get-LinkAddress $sheet $cRef

    $aText = $sheet.Range("A" + $row).Value()
    $cText = $sheet.Range("C" + $row).Value()

    $sheet.Range("B" + $row).Value = "Handed back: in sync with en-US"

    $sheet.Hyperlinks.Add($sheet.Range("E" + $row), $aAddr, "", "", $aText)
    $sheet.Range("E" + $row).Font.Underline = 2
    $sheet.Range("E" + $row).Font.Color = 15570276

    $sheet.Hyperlinks.Add($sheet.Range("F" + $row), $cAddr, "", "", $cText)
    $sheet.Range("F" + $row).Font.Underline = 2
    $sheet.Range("F" + $row).Font.Color = 15570276

    $sheet.Range("G" + $row).Value = $handbackDatetime
}

$zhcn = $wb.Worksheets.Item("zh-cn")
Set-HandbackRow $zhcn 2 "2016-01-26 06:24:26"
Set-HandbackRow $zhcn 3 "2016-01-26 06:24:26"

$dede = $wb.Worksheets.Item("de-de")
Set-HandbackRow $dede 2 "2016-01-26 06:24:50"
Set-HandbackRow $dede 3 "2016-01-26 06:24:50"

# The "Status" text is a shared string also surfaced (via the same value)
# on the "Overview" summary sheet's B/C columns for these two rows - update
# those too so the whole workbook reflects the new status consistently.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"
